# chore: update Sheets via scheduled runner
#
# This script refreshes computed market/profit figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -> columns H:N)
# for a handful of leve rows across several job sheets, mirroring a
# scheduled data-refresh run.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")

# Row 76
$ws.Range("H76").Value = 3375.125
$ws.Range("I76").Value = 2501
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 2501
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -2186
$ws.Range("N76").Value = -4130

# Row 79
$ws.Range("H79").Value = 3375.125
$ws.Range("I79").Value = 2501
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 2501
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -1409
$ws.Range("N79").Value = -5684

# Row 98
$ws.Range("H98").Value = 1389.0435
$ws.Range("I98").Value = 1188
$ws.Range("J98").Value = 3500
$ws.Range("K98").Value = 1188
$ws.Range("L98").Value = 3500
$ws.Range("M98").Value = 310
$ws.Range("N98").Value = -6496

# Row 113
$ws.Range("H113").Value = 1869.2307
$ws.Range("I113").Value = 1575
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1575
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1679
$ws.Range("N113").Value = -8508

# Row 116
$ws.Range("H116").Value = 3207723.8
$ws.Range("I116").Value = 25643358
$ws.Range("J116").Value = 2633.3333
$ws.Range("K116").Value = 25643358
$ws.Range("L116").Value = 2633.3333
$ws.Range("M116").Value = -25639916
$ws.Range("N116").Value = -9517.3333

# Row 122
$ws.Range("H122").Value = 1389.0435
$ws.Range("I122").Value = 1188
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 3564
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -1114
$ws.Range("N122").Value = -15400

# Row 125
$ws.Range("H125").Value = 715.1
$ws.Range("I125").Value = 490.5
$ws.Range("J125").Value = 1052
$ws.Range("K125").Value = 4414.5
$ws.Range("L125").Value = 9468
$ws.Range("M125").Value = -1954.5
$ws.Range("N125").Value = -14388

# Row 129
$ws.Range("H129").Value = 1415.3549
$ws.Range("I129").Value = 478
$ws.Range("J129").Value = 2092.3333
$ws.Range("K129").Value = 1434
$ws.Range("L129").Value = 6276.999899999999
$ws.Range("M129").Value = 3566
$ws.Range("N129").Value = -16276.9999

# Row 132
$ws.Range("H132").Value = 4149.972
$ws.Range("I132").Value = 4046.76
$ws.Range("J132").Value = 4384.5454
$ws.Range("K132").Value = 12140.28
$ws.Range("L132").Value = 13153.6362
$ws.Range("M132").Value = -9610.280000000001
$ws.Range("N132").Value = -18213.6362

$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 14709701
$ws.Range("I32").Value = 16131900
$ws.Range("J32").Value = 13647.167
$ws.Range("K32").Value = 16131900
$ws.Range("L32").Value = 13647.167
$ws.Range("M32").Value = -16131613
$ws.Range("N32").Value = -14221.167

# Row 61
$ws.Range("H61").Value = 1791.762
$ws.Range("I61").Value = 1141.8572
$ws.Range("J61").Value = 3091.5715
$ws.Range("K61").Value = 1141.8572
$ws.Range("L61").Value = 3091.5715
$ws.Range("M61").Value = -929.8571999999999
$ws.Range("N61").Value = -3515.5715

# Row 122
$ws.Range("H122").Value = 1123.44
$ws.Range("I122").Value = 1008.6667
$ws.Range("J122").Value = 1726
$ws.Range("K122").Value = 3026.0001
$ws.Range("L122").Value = 5178
$ws.Range("M122").Value = -576.0001000000002
$ws.Range("N122").Value = -10078

# Row 132
$ws.Range("H132").Value = 1489.6296
$ws.Range("I132").Value = 1341.8529
$ws.Range("J132").Value = 1740.85
$ws.Range("K132").Value = 4025.5587
$ws.Range("L132").Value = 5222.549999999999
$ws.Range("M132").Value = -1495.5587
$ws.Range("N132").Value = -10282.55

# Row 136
$ws.Range("H136").Value = 1791.762
$ws.Range("I136").Value = 1141.8572
$ws.Range("J136").Value = 3091.5715
$ws.Range("K136").Value = 3425.5716
$ws.Range("L136").Value = 9274.7145
$ws.Range("M136").Value = -875.5715999999998
$ws.Range("N136").Value = -14374.7145

$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 5959.625
$ws.Range("I20").Value = 6899
$ws.Range("J20").Value = 5646.5
$ws.Range("K20").Value = 6899
$ws.Range("L20").Value = 5646.5
$ws.Range("M20").Value = -6652
$ws.Range("N20").Value = -6140.5

# Row 134
$ws.Range("H134").Value = 1770
$ws.Range("I134").Value = 1417.579
$ws.Range("J134").Value = 3109.2
$ws.Range("K134").Value = 4252.737
$ws.Range("L134").Value = 9327.599999999999
$ws.Range("M134").Value = -1717.737
$ws.Range("N134").Value = -14397.6

$ws = $wb.Worksheets.Item("CRP")

# Row 99
$ws.Range("H99").Value = 2779.8
$ws.Range("I99").Value = 3133
$ws.Range("J99").Value = 2250
$ws.Range("K99").Value = 3133
$ws.Range("L99").Value = 2250
$ws.Range("M99").Value = -1635
$ws.Range("N99").Value = -5246

# Row 126
$ws.Range("H126").Value = 2779.8
$ws.Range("I126").Value = 3133
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 9399
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -6929
$ws.Range("N126").Value = -11690

# Row 132
$ws.Range("H132").Value = 1808.9
$ws.Range("I132").Value = 1350.2
$ws.Range("J132").Value = 4102.4
$ws.Range("K132").Value = 4050.6
$ws.Range("L132").Value = 12307.2
$ws.Range("M132").Value = -1520.6
$ws.Range("N132").Value = -17367.2

# Row 134
$ws.Range("H134").Value = 1855.125
$ws.Range("I134").Value = 1146.8372
$ws.Range("J134").Value = 7946.4
$ws.Range("K134").Value = 3440.5116
$ws.Range("L134").Value = 23839.2
$ws.Range("M134").Value = -905.5115999999998
$ws.Range("N134").Value = -28909.2

$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 5129.159
$ws.Range("I70").Value = 5158.5386
$ws.Range("J70").Value = 4900
$ws.Range("K70").Value = 5158.5386
$ws.Range("L70").Value = 4900
$ws.Range("M70").Value = -4888.5386
$ws.Range("N70").Value = -5440

# Row 73
$ws.Range("H73").Value = 5129.159
$ws.Range("I73").Value = 5158.5386
$ws.Range("J73").Value = 4900
$ws.Range("K73").Value = 5158.5386
$ws.Range("L73").Value = 4900
$ws.Range("M73").Value = -4222.5386
$ws.Range("N73").Value = -6772

# Row 105
$ws.Range("H105").Value = 90000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 90000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 90000
$ws.Range("N105").Value = -96988

# Row 132
$ws.Range("H132").Value = 2896.9387
$ws.Range("I132").Value = 2700.4167
$ws.Range("J132").Value = 3441.1538
$ws.Range("K132").Value = 8101.250100000001
$ws.Range("L132").Value = 10323.4614
$ws.Range("M132").Value = -5571.250100000001
$ws.Range("N132").Value = -9517.3333

$ws = $wb.Worksheets.Item("LTW")

# Row 40
$ws.Range("H40").Value = 3698.3948
$ws.Range("I40").Value = 2660.5386
$ws.Range("J40").Value = 5947.0835
$ws.Range("K40").Value = 2660.5386
$ws.Range("L40").Value = 5947.0835
$ws.Range("M40").Value = -2524.5386
$ws.Range("N40").Value = -6219.0835

# Row 132
$ws.Range("H132").Value = 1059.2211
$ws.Range("I132").Value = 1056.5222
$ws.Range("J132").Value = 1107.8
$ws.Range("K132").Value = 3169.5666
$ws.Range("L132").Value = 3323.4
$ws.Range("M132").Value = -1520.6
